$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("ECs","Cd274","Cd80","ECs",3,1,17.56563533333333,52.696906,0.8967254513503751,0.8967254513503751,3,1,30.40041266666667,91.201238,0.8640973522824783,0.8640973522824784,534.0025628855142,4806.023065969628,0.7748580882361694,0.7748580882361695)
    ,@("ECs","Cd274","Cd80","FAPs",3,1,17.56563533333333,52.696906,0.8967254513503751,0.8967254513503751,3,1,4.356394000000001,13.069182,0.1238255730991045,0.1238255730991045,76.52282837232134,688.7054553508921,0.1110375429260133,0.1110375429260133)
    ,@("ECs","Cd274","Cd80","sCs",3,1,17.56563533333333,52.696906,0.8967254513503751,0.8967254513503751,3,1,0.424892,1.274676,0.01207707461841714,0.01207707461841714,7.463497928050668,67.171481352456,0.01082982018819227,0.01082982018819227)
    ,@("FAPs","Cd274","Cd80","ECs",2,0.6666666666666666,1.33491,4.00473,0.06814713783739766,0.06814713783739766,3,1,30.40041266666667,91.201238,0.8640973522824783,0.8640973522824784,40.58181487286,365.23633385574,0.05888576137092441,0.05888576137092442)
    ,@("FAPs","Cd274","Cd80","FAPs",2,0.6666666666666666,1.33491,4.00473,0.06814713783739766,0.06814713783739766,3,1,4.356394000000001,13.069182,0.1238255730991045,0.1238255730991045,5.815393914540001,52.33854523086001,0.008438358397779434,0.008438358397779434)
    ,@("FAPs","Cd274","Cd80","sCs",2,0.6666666666666666,1.33491,4.00473,0.06814713783739766,0.06814713783739766,3,1,0.424892,1.274676,0.01207707461841714,0.01207707461841714,0.56719257972,5.104733217480001,0.0008230180686938094,0.0008230180686938094)
    ,@("sCs","Cd274","Cd80","ECs",3,1,0.6880983333333334,2.064295,0.03512741081222724,0.03512741081222724,3,1,30.40041266666667,91.201238,0.8640973522824783,0.8640973522824784,20.91847328857889,188.26625959721,0.03035350267538446,0.03035350267538446)
    ,@("sCs","Cd274","Cd80","FAPs",3,1,0.6880983333333334,2.064295,0.03512741081222724,0.03512741081222724,3,1,4.356394000000001,13.069182,0.1238255730991045,0.1238255730991045,2.997627450743334,26.97864705669,0.004349671775311718,0.004349671775311718)
    ,@("sCs","Cd274","Cd80","sCs",3,1,0.6880983333333334,2.064295,0.03512741081222724,0.03512741081222724,3,1,0.424892,1.274676,0.01207707461841714,0.01207707461841714,0.2923674770466667,2.63130729342,0.0004242363615310613,0.0004242363615310613)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
